$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 1. Update the "Date Colours" sheet colour codes (B2:B4) to new shades ---
$ws2.Range("B2").Value = "#d9d6ec"
$ws2.Range("B3").Value = "#b4afd9"
$ws2.Range("B4").Value = "#8e8ac5"

# --- 2. Rename node labels on Sheet1 (whole-cell, case-sensitive replace) ---
$rng1 = $ws1.UsedRange
$rng1.Replace("M7", "M7 m50", -4143, 1, $true, $false, $true) | Out-Null
$rng1.Replace("M8", "M8 m60", -4143, 1, $true, $false, $true) | Out-Null
$rng1.Replace("M9", "M9 9yo", -4143, 1, $true, $false, $true) | Out-Null
$rng1.Replace("M1", "M1 30", -4143, 1, $true, $false, $true) | Out-Null
$rng1.Replace("M2", "M2 30", -4143, 1, $true, $false, $true) | Out-Null
$rng1.Replace("M3", "M3 child", -4143, 1, $true, $false, $true) | Out-Null

# --- 3. Append two new outbreak-path rows to the Sheet1 table ---
$lo = $ws1.ListObjects.Item(1)

$lo.ListRows.Add() | Out-Null
$ws1.Range("A16").NumberFormat = "d-mmm"
$ws1.Range("A16").Value = 44393
$ws1.Range("B16").Value = "M4 m60"
$ws1.Range("C16").Value = "M10 11yo"
$ws1.Range("D16").Value = "Maribyrnong"
$ws1.Range("E16").Value = "MCG"
$ws1.Range("F16").Value = "MCG"
$ws1.Range("G16").Value = "Delta (B.1.617.2)"
$ws1.Range("H16").Value = "Not Isolated"

$lo.ListRows.Add() | Out-Null
$ws1.Range("A17").NumberFormat = "d-mmm"
$ws1.Range("A17").Value = 44393
$ws1.Range("B17").Value = "M4 m60"
$ws1.Range("C17").Value = "M11 Adult"
$ws1.Range("D17").Value = "Maribyrnong"
$ws1.Range("E17").Value = "MCG"
$ws1.Range("F17").Value = "MCG"
$ws1.Range("G17").Value = "Delta (B.1.617.2)"
$ws1.Range("H17").Value = "Not Isolated"

# --- 4. New colour entry on the "Date Colours" sheet for the newly added date ---
$ws2.Range("B5").Value = "#6666b2"

# --- 5. Update selections / active sheet to match the saved view state ---
$ws2.Range("F2:I2").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("H17").Select() | Out-Null
